$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 3:4 (this shifts the current rows 3:8 down to 5:10).
$ws.Rows("3:4").Insert()

# The rows that used to be 7:8 are now at 9:10. Move (copy+delete) that
# content into the freshly inserted rows 3:4, then remove the now-duplicated
# rows 9:10, which restores the original row count.
$ws.Range("A9:R10").Copy()
$ws.Range("A3").PasteSpecial(-4163)

# Re-apply the couple of cells whose style differs from the sheet/column
# default (style 1), since a bulk "paste values" already yields style 1.
$ws.Range("N9").Copy()
$ws.Range("N3").PasteSpecial(-4122)
$ws.Range("R9").Copy()
$ws.Range("R3").PasteSpecial(-4122)

$ws.Rows("9:10").Delete()

$ws.Range("A3:XFD4").Select()
